$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 272.57144
$ws.Range("I38").Value = 201.33333
$ws.Range("J38").Value = 700
$ws.Range("K38").Value = 603.99999
$ws.Range("L38").Value = 2100
$ws.Range("M38").Value = -231.99999
$ws.Range("N38").Value = -2844
$ws.Range("H43").Value = 826
$ws.Range("I43").Value = 750
$ws.Range("J43").Value = 845
$ws.Range("K43").Value = 750
$ws.Range("L43").Value = 845
$ws.Range("M43").Value = -681
$ws.Range("N43").Value = -983
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9508
$ws.Range("H132").Value = 286371.44
$ws.Range("I132").Value = 338467.7
$ws.Range("K132").Value = 1015403.1
$ws.Range("M132").Value = -1012873.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 66200
$ws.Range("J80").Value = 49800
$ws.Range("L80").Value = 49800
$ws.Range("N80").Value = -51796
$ws.Range("H83").Value = 66200
$ws.Range("J83").Value = 49800
$ws.Range("L83").Value = 149400
$ws.Range("N83").Value = -159384
$ws.Range("H122").Value = 1122.7368
$ws.Range("I122").Value = 1000.7692
$ws.Range("J122").Value = 1387
$ws.Range("K122").Value = 3002.3076
$ws.Range("L122").Value = 4161
$ws.Range("M122").Value = -552.3076000000001
$ws.Range("N122").Value = -9061
$ws.Range("H132").Value = 2272
$ws.Range("I132").Value = 1973.1111
$ws.Range("J132").Value = 4962
$ws.Range("K132").Value = 5919.3333
$ws.Range("L132").Value = 14886
$ws.Range("M132").Value = -3389.3333
$ws.Range("N132").Value = -19946

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100630
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102184
$ws.Range("H107").Value = 630.0476
$ws.Range("I107").Value = 623.2308
$ws.Range("K107").Value = 623.2308
$ws.Range("M107").Value = 1296.7692
$ws.Range("H132").Value = 19800
$ws.Range("J132").Value = 19800
$ws.Range("L132").Value = 19800
$ws.Range("N132").Value = -29920
$ws.Range("H138").Value = 49466.668
$ws.Range("J138").Value = 49466.668
$ws.Range("L138").Value = 49466.668
$ws.Range("N138").Value = -59746.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 166.16667
$ws.Range("I22").Value = 172
$ws.Range("J22").Value = 102
$ws.Range("K22").Value = 172
$ws.Range("L22").Value = 102
$ws.Range("M22").Value = 178
$ws.Range("N22").Value = -802
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H134").Value = 4485.2856
$ws.Range("I134").Value = 2504.8
$ws.Range("J134").Value = 5585.5557
$ws.Range("K134").Value = 7514.400000000001
$ws.Range("L134").Value = 16756.6671
$ws.Range("M134").Value = -4979.400000000001
$ws.Range("N134").Value = -21826.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 5013
$ws.Range("I75").Value = 5013
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 15039
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -14041
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 5013
$ws.Range("I78").Value = 5013
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 45117
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -40125
$ws.Range("N78").ClearContents()
$ws.Range("H113").Value = 872.0303
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 899.23334
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2697.70002
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -7037.70002
$ws.Range("H121").Value = 730
$ws.Range("I121").Value = 190
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 570
$ws.Range("L121").Value = 3000
$ws.Range("M121").Value = 740
$ws.Range("N121").Value = -5620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 12800
$ws.Range("J17").Value = 3500
$ws.Range("L17").Value = 3500
$ws.Range("N17").Value = -3836
$ws.Range("H102").Value = 2014.2858
$ws.Range("I102").Value = 1620
$ws.Range("K102").Value = 1620
$ws.Range("M102").Value = 2
$ws.Range("H107").Value = 338
$ws.Range("J107").Value = 351
$ws.Range("L107").Value = 351
$ws.Range("N107").Value = -4191
$ws.Range("H132").Value = 2475.3901
$ws.Range("I132").Value = 1986.6666
$ws.Range("K132").Value = 5959.9998
$ws.Range("M132").Value = -3429.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 27500
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5280
$ws.Range("H22").Value = 17366.834
$ws.Range("I22").Value = 733
$ws.Range("J22").Value = 34000.668
$ws.Range("K22").Value = 733
$ws.Range("L22").Value = 34000.668
$ws.Range("M22").Value = -438
$ws.Range("N22").Value = -34590.668
$ws.Range("H27").Value = 17366.834
$ws.Range("I27").Value = 733
$ws.Range("J27").Value = 34000.668
$ws.Range("K27").Value = 733
$ws.Range("L27").Value = 34000.668
$ws.Range("M27").Value = -626
$ws.Range("N27").Value = -34214.668
$ws.Range("H55").Value = 236.125
$ws.Range("I55").Value = 210.45833
$ws.Range("J55").Value = 313.125
$ws.Range("K55").Value = 210.45833
$ws.Range("L55").Value = 313.125
$ws.Range("M55").Value = -37.45832999999999
$ws.Range("N55").Value = -659.125
$ws.Range("H93").Value = 480.5625
$ws.Range("I93").Value = 507
$ws.Range("K93").Value = 507
$ws.Range("M93").Value = 741
$ws.Range("H132").Value = 3113.611
$ws.Range("I132").Value = 1913.3636
$ws.Range("K132").Value = 5740.0908
$ws.Range("M132").Value = -3210.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12823781
$ws.Range("I132").Value = 21743086
$ws.Range("J132").Value = 2280.1875
$ws.Range("K132").Value = 65229258
$ws.Range("L132").Value = 6840.5625
$ws.Range("M132").Value = -65226728
$ws.Range("N132").Value = -11900.5625
